# "Generate Report for Handoff" - refresh the localization-status report:
#   - Overview / zh-cn / de-de sheets: status text flips from the old
#     "handed back" wording to "Ready for handoff", and the associated
#     timestamps move forward a bit (a new handoff was just generated).
#   - The "Status" columns get narrower now that the new text is shorter
#     than the old one, so re-apply their (auto-fit) widths. Excel's
#     ColumnWidth setter only lands on whole-pixel steps, so use the
#     input that rounds to the pixel closest to the target character
#     width (~17.22 chars).
$targetStatusWidth = 16.3333333333333

$wb = $excel.ActiveWorkbook

# --- Overview sheet ------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-08-13 13:15:11"
$ov.Range("E1:F1").EntireColumn.ColumnWidth = $targetStatusWidth

# --- zh-cn sheet -----------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("H2").Value = "2016-08-13 13:15:00"
$zh.Range("C1").EntireColumn.ColumnWidth = $targetStatusWidth

# --- de-de sheet -------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Ready for handoff"
$de.Range("H2").Value = "2016-08-13 13:15:11"
$de.Range("C1").EntireColumn.ColumnWidth = $targetStatusWidth
